$d = $word.ActiveDocument

# Fix the misspelling "Manageement" -> "Management" in the
# "Account Manageement (H4)" heading line (spell-check correction).
$d.Content.Find.Execute("Account Manageement (H4)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Account Management (H4)", 2)

# Locate the corrected phrase so we can reproduce the resulting run split
# ("Account Manag" / "ement (H4)") that Word leaves behind after the
# spelling correction.
$found = $d.Content
$found.Find.Execute("Account Management (H4)", $true, $false, $false, $false, $false,
                     $true, 1, $false, "", 0)
$phraseStart = $found.Start

# "Account Manag" is 13 characters long; nudging the character formatting
# on the remainder forces Word to split the text into two runs at that
# boundary, matching the post-edit document structure.
$splitPoint = $phraseStart + 13
$tail = $d.Range($splitPoint, $phraseStart + 23)
$tail.Font.Bold = 1
$tail.Font.Bold = 0
